$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 8: AUTRE category with no Value
$ws.Range("A8").Value = "AUTRE"

# Update B5 (LOISIRS) value: append ";micromania"
$ws.Range("B5").Value = "Easy cash;Call of duty;Nintendo;Instant gaming;leboncoin;fnac;figurines;micromania"

# Update selection to match the diff (activeCell B5)
$ws.Range("B5").Select()
